# Refatorando o consolidador para modelo ETL
# Replace the absenteeism sample rows (2-11) with the newly generated dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=2765;  B="Joaquim da Conceição";    C="Marketing";   D="Viagem de negócios"; E=2; F=45085; G=8632.370000000001},
    @{Row=3;  A=40654; B="Dra. Nina da Mata";        C="P&D";         D="Doença";             E=4; F=45101; G=3942.77},
    @{Row=4;  A=9455;  B="Camila Correia";           C="Engenharia";  D="Viagem de negócios"; E=4; F=45096; G=10267.75},
    @{Row=5;  A=96023; B="Yago Monteiro";            C="Marketing";   D="Consulta médica";    E=4; F=45101; G=6568.45},
    @{Row=6;  A=62544; B="Maria Vitória da Rocha";   C="Jurídico";    D="Outros";             E=6; F=45089; G=6758.57},
    @{Row=7;  A=52293; B="Sr. Fernando Moreira";     C="Marketing";   D="Problemas pessoais"; E=8; F=45084; G=6907.86},
    @{Row=8;  A=91110; B="Srta. Bruna da Costa";     C="TI";          D="Doença";             E=5; F=45087; G=2884.04},
    @{Row=9;  A=25867; B="Kaique Pinto";             C="Operações";   D="Problemas pessoais"; E=4; F=45094; G=3420.86},
    @{Row=10; A=72185; B="Francisco Gomes";          C="TI";          D="Doença";             E=6; F=45096; G=6312},
    @{Row=11; A=76451; B="Davi Luiz Moreira";        C="Operações";   D="Viagem de negócios"; E=2; F=45079; G=5545.42}
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.A
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
}
